$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated header labels: fix "chai" -> "Chai" casing, reword two labels, and
# make every header run bold (matching the already-bold "Date" header run).
$headers = [ordered]@{
    "B1" = "Ventas totales de Chai (unidades)"
    "C1" = "Ventas de Chai artesanal (unidades)"
    "D1" = "Ventas predefinidas de Chai (unidades)"
    "E1" = "Interacción de redes sociales (visualizaciones)"
    "F1" = "Búsquedas en línea de Chai"
}

$white = 16777215

foreach ($addr in $headers.Keys) {
    $rng = $ws.Range($addr)

    # Update the text. (This replaces the cell's shared-string content.)
    $rng.Value = $headers[$addr]

    $len = $rng.Text.Length

    # Restore the header's rich-text look (white, bold, Calibri 11) across the
    # whole string. Splitting into two adjoining character runs avoids a
    # "whole string" fast path and reliably yields per-run formatting that
    # Excel then stores as a single merged <r> run, same shape as "Date".
    $rng.Characters(1, $len - 1).Font.Color = $white
    $rng.Characters($len, 1).Font.Color = $white

    $rng.Characters(1, $len - 1).Font.Name = "Calibri"
    $rng.Characters($len, 1).Font.Name = "Calibri"

    $rng.Characters(1, $len - 1).Font.Size = 11
    $rng.Characters($len, 1).Font.Size = 11

    $rng.Characters(1, $len - 1).Font.Bold = $true
    $rng.Characters($len, 1).Font.Bold = $true
}
